# "minor updates to lectures 3 and 4"
#
# Slide 8 ("What are raw counts") - the title is reworded into a question:
#   "What are raw counts"  ->  "What are raw " + "counts?"
# In the saved OOXML this shows up as the single run being split into two
# runs (both carrying the same rPr), which is exactly what happens when you
# position the cursor in the middle of an existing run and type/insert new
# text in PowerPoint. We reproduce that here by grabbing the "counts"
# sub-range of the title and inserting "?" immediately after it - this
# naturally splits "What are raw counts" into "What are raw " + "counts?"
# without disturbing the run's formatting (rPr stays identical on both
# halves, matching the diff).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$titleShape = $s.Shapes.Item(1)
$tr = $titleShape.TextFrame.TextRange

$fullText = $tr.Text
$idx = $fullText.IndexOf("counts")
if ($idx -ge 0) {
    $startPos = $idx + 1            # TextRange.Characters is 1-based
    $len = "counts".Length
    $wordRange = $tr.Characters($startPos, $len)
    $wordRange.InsertAfter("?")
} else {
    # Fallback in case the text ever differs from what we expect.
    $tr.Text = "What are raw counts?"
}
